# Adds a new "2022-Q3" sheet (with its fund-holding detail data) right after
# "总计" and before "2022-Q2", and updates the "总计" (totals) sheet with the
# new quarter's row, shifting the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: force a numeric-looking string to be stored as TEXT (matching the
# workbook's convention of keeping these columns as inline strings) and then
# strip the left-over "@" number-format styling so the cell ends up with the
# plain/default style - exactly like the pre-existing sibling cells.
# ---------------------------------------------------------------------------
# NOTE: always invoke with POSITIONAL arguments - named arguments
# (-Cell/-Value/-Scratch) do not bind correctly in this PowerShell host.
function Set-TextValue {
    param($Cell, $Value, $Scratch)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Scratch.Copy()
    $Cell.PasteSpecial(-4122)
}

# ===========================================================================
# 1. Insert the new "2022-Q3" worksheet before the current "2022-Q2" sheet.
# ===========================================================================
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# A blank, never-touched cell used purely as a "default style" format donor.
$scratch = $newSheet.Cells.Item(200, 60)

# Copy the header-row formatting (bold/border/centered = style used by every
# other sheet's row 1) from the sheet that is about to become "2022-Q2".
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Copy the column-A "index" cell formatting (bold/border/centered, style
# shared with every other sheet's column A) as well.
$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A3:A16").Value = $newSheet.Range("A2").Value
$newSheet.Range("A2").Copy()
$newSheet.Range("A2:A16").PasteSpecial(-4122)

# Header row text.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Fund-holding detail rows for 2022-Q3.
$data = @(
    @(0, "519087", "新华优选分红混合", "10.82", "89.55", "5.04", "0.5453", 6),
    @(1, "160211", "国泰中小盘成长混合（LOF）", "6.30", "87.18", "6.94", "0.4372", 1),
    @(2, "001040", "新华策略精选股票", "6.78", "94.54", "5.04", "0.3417", 6),
    @(3, "003231", "创金合信医疗保健行业股票C", "6.26", "94.64", "5.33", "0.3337", 10),
    @(4, "519156", "新华行业轮换灵活配置混合A", "5.70", "94.21", "5.26", "0.2998", 6),
    @(5, "003230", "创金合信医疗保健行业股票A", "3.66", "94.64", "5.33", "0.1951", 10),
    @(6, "506009", "国泰科创板两年定期开放混合", "2.05", "85.80", "6.60", "0.1353", 1),
    @(7, "014126", "华夏中证1000指数增强C", "8.78", "89.62", "0.84", "0.0738", 2),
    @(8, "001294", "新华战略新兴产业灵活配置混合", "0.99", "93.49", "5.55", "0.0549", 3),
    @(9, "011457", "新华行业龙头主题股票", "0.97", "94.28", "5.41", "0.0525", 6),
    @(10, "005520", "国投瑞银创新医疗混合", "0.41", "93.90", "3.15", "0.0129", 10),
    @(11, "005997", "天弘裕利灵活配置混合C", "0.50", "44.05", "1.94", "0.0097", 1),
    @(12, "014125", "华夏中证1000指数增强A", "0.97", "89.62", "0.84", "0.0081", 2),
    @(13, "519157", "新华行业轮换灵活配置混合C", "0.06", "94.21", "5.26", "0.0032", 6),
    @(14, "002388", "天弘裕利灵活配置混合A", "0.10", "44.05", "1.94", "0.0019", 1)
)

foreach ($row in $data) {
    $r = [int]$row[0] + 2
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $newSheet.Cells.Item($r, 2) $row[1] $scratch
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[3] $scratch
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[4] $scratch
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[5] $scratch
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[6] $scratch
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# ===========================================================================
# 2. Update the "总计" (totals) summary sheet: insert the 2022-Q3 row at the
#    top of the data and shift every other quarter down by one row.
# ===========================================================================
$total = $wb.Worksheets.Item("总计")

$totalsData = @(
    @(0, "2022-Q3", 15, 2.51),
    @(1, "2022-Q2", 7, 0.89),
    @(2, "2022-Q1", 3, 0.25),
    @(3, "2021-Q4", 1, 0.16),
    @(4, "2021-Q2", 15, 0.58),
    @(5, "2021-Q1", 3, 0.03),
    @(6, "2020-Q4", 3, 0.04)
)

# Make sure row 8 (brand new) has the same "index column" style as the rows
# above it before we populate it.
$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial(-4122)

foreach ($row in $totalsData) {
    $r = [int]$row[0] + 2
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}
